$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.263.49"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "1.916.26"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").NumberFormat = "general"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7384"
$ws.Range("D5").NumberFormat = "general"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.25"
$ws.Range("D6").NumberFormat = "general"
$ws.Range("E6").Value = "  -3.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.006"
$ws.Range("D7").NumberFormat = "general"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3121"
$ws.Range("D8").NumberFormat = "general"
$ws.Range("E8").Value = "  -2.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.09"
$ws.Range("D9").NumberFormat = "general"
$ws.Range("E9").Value = "  -2.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06926"
$ws.Range("D10").NumberFormat = "general"
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07997"
$ws.Range("D11").NumberFormat = "general"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7659"
$ws.Range("D12").NumberFormat = "general"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").Value = "1.932.99"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.280"
$ws.Range("D14").NumberFormat = "general"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.14"
$ws.Range("D15").NumberFormat = "general"
$ws.Range("E15").Value = "  -3.74%  "
$ws.Range("D16").Value = "30.295.73"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.12"
$ws.Range("D17").NumberFormat = "general"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "246.25"
$ws.Range("D18").NumberFormat = "general"
$ws.Range("E18").Value = "  -3.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.824"
$ws.Range("D19").NumberFormat = "general"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007820"
$ws.Range("D20").NumberFormat = "general"
$ws.Range("E20").Value = "  -2.79%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.186.94"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.006"
$ws.Range("D22").NumberFormat = "general"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("D23").NumberFormat = "general"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.578"
$ws.Range("D24").NumberFormat = "general"
$ws.Range("E24").Value = "  -3.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.347"
$ws.Range("D25").NumberFormat = "general"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.07"
$ws.Range("D26").NumberFormat = "general"
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.83"
$ws.Range("D27").NumberFormat = "general"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1274"
$ws.Range("D28").NumberFormat = "general"
$ws.Range("E28").Value = "  -3.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.133"
$ws.Range("D29").NumberFormat = "general"
$ws.Range("E29").Value = "  -7.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.357"
$ws.Range("D30").NumberFormat = "general"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.538"
$ws.Range("D31").NumberFormat = "general"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("E32").Value = "  -2.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.045"
$ws.Range("D33").NumberFormat = "general"
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05148"
$ws.Range("D34").NumberFormat = "general"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.293"
$ws.Range("D35").NumberFormat = "general"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7417"
$ws.Range("D36").NumberFormat = "general"
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.776"
$ws.Range("D37").NumberFormat = "general"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01931"
$ws.Range("D38").NumberFormat = "general"
$ws.Range("E38").Value = "  -2.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.773"
$ws.Range("D39").NumberFormat = "general"
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.401"
$ws.Range("D40").NumberFormat = "general"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.61"
$ws.Range("D41").NumberFormat = "general"
$ws.Range("E41").Value = "  -3.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4436"
$ws.Range("D42").NumberFormat = "general"
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.925"
$ws.Range("D43").NumberFormat = "general"
$ws.Range("E43").Value = "  -3.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.005"
$ws.Range("D44").NumberFormat = "general"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8360"
$ws.Range("D45").NumberFormat = "general"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.86"
$ws.Range("D46").NumberFormat = "general"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.548"
$ws.Range("D47").NumberFormat = "general"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.765"
$ws.Range("D48").NumberFormat = "general"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.93"
$ws.Range("D49").NumberFormat = "general"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "952.27"
$ws.Range("D50").NumberFormat = "general"
$ws.Range("E50").Value = "  -2.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1201"
$ws.Range("D51").NumberFormat = "general"
$ws.Range("E51").Value = "  +4.38%  "
